$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Murphy")

# Row 26 (S/N 25): "Front-end: Documentation and refactoring"
# Replace the old comment text with a new one describing delayed work
$ws.Range("H26").Value = "Delayed as I worked on integrating with the backend API first`nRedux implemented into the general part of the app, currently working on using it with the three.js scene"

# Row 27 (S/N 26): "Front-end: Integrating with backend API"
# Duration changes from "?" to 6, Finish/Actual-Finish change from "?" to real dates
$ws.Range("D27").Value = 6

# Copy the date number format from column E (Start) so the new date cells
# share the same style index instead of minting a brand-new number format
$ws.Range("E27").Copy()
$ws.Range("F27").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G27").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("F27").Value2 = 44321
$ws.Range("G27").Value2 = 44321
$ws.Range("H27").Value = "Current data used from API (connected locally with docker)`n- list of areas`n- list of blocks (per area)`n- list of lights (per block)`n- status of light (on/off/dimmed)`n- fault status of light`nTime estimate is for the current data available in the backend API"

# Update row heights to match re-wrapped comment text
$ws.Rows.Item(26).RowHeight = 75
$ws.Rows.Item(27).RowHeight = 120

# Update view / selection state to match where the user ended up working
$ws.Range("H28").Select()
